$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text so numeric-looking strings (e.g. "1.00")
# are preserved verbatim as text, matching the original inlineStr cells,
# then restore the original (General/Normal) styling afterwards.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "90.736.52"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.119.35"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "238.98"
$ws.Range("E5").Value = "  +9.84%  "
$ws.Range("D6").Value = "631.31"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.361"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.117.89"
$ws.Range("D11").Value = "0.715"
$ws.Range("E11").Value = "  -5.46%  "
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "36.68"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "5.50"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "90.599.44"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "3.693.40"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "3.161.71"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "3.81"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "14.17"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "0.0000208"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").Value = "443.10"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "5.56"
$ws.Range("E23").Value = "  +6.54%  "
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "5.72"
$ws.Range("E25").Value = "  -8.47%  "
$ws.Range("D26").Value = "12.61"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "88.36"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "3.304.88"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "9.56"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").Value = "0.160"
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "26.35"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "0.194"
$ws.Range("E33").Value = "  +21.97%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.889"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "3.80"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("D36").Value = "509.48"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").Value = "0.148"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "7.17"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "1.29"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "0.412"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "22.17"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "0.0844"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("E45").Value = "  +48.40%  "
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "151.37"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "0.684"
$ws.Range("E48").Value = "  +6.48%  "
$ws.Range("D49").Value = "45.03"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "4.47"
$ws.Range("E51").Value = "  +1.77%  "

$dRange.Style = "Normal"
